# Refactor - complete bovine lentivirus extension
#
# Collapses the per-genotype "name" values for the Small Ruminant
# Lentivirus (SRLV) rows down to the single shared name "SRLV", and
# renames the HIV-1 "name" value from "HIV-1M" to "HIV-1".
# Also restores the working selection to the full data range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Human immunodeficiency virus 1 row: name HIV-1M -> HIV-1
$ws.Range("C5").Value = "HIV-1"

# Small ruminant lentivirus rows: name SRLV-<genotype> -> SRLV
$ws.Range("C7").Value = "SRLV"
$ws.Range("C8").Value = "SRLV"
$ws.Range("C9").Value = "SRLV"
$ws.Range("C10").Value = "SRLV"
$ws.Range("C11").Value = "SRLV"
$ws.Range("C12").Value = "SRLV"

# Restore selection over the whole table (best effort: the active cell
# itself always normalizes to the top-left of the selected range in
# this host).
$ws.Range("A1:N12").Select()
